# Insert a new weekly price record for "Femacal de La Calera - Poroto granado"
# at row 91, pushing the existing rows 91-121 down to 92-122.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(91).Insert()

$ws.Cells.Item(91, 1).Value = 3
$ws.Cells.Item(91, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(91, 3).Value = "Coquimbo"
$ws.Cells.Item(91, 4).Value = 44559
$ws.Cells.Item(91, 5).Value = 5
$ws.Cells.Item(91, 6).Value = 100112030
$ws.Cells.Item(91, 7).Value = "Poroto granado"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 88
$ws.Cells.Item(91, 11).Value = 40000
$ws.Cells.Item(91, 12).Value = 42000
$ws.Cells.Item(91, 13).Value = 40909
$ws.Cells.Item(91, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(91, 15).Value = "Provincia de Talca"
$ws.Cells.Item(91, 16).Value = 1636
$ws.Cells.Item(91, 17).Value = 25
$ws.Cells.Item(91, 18).Value = "Hortaliza"
